# L2.docx edits ("fargade saker i l2")
#
# 1) Merge the three footer runs that make up "Rasmus Tilljander -
#    rati10@student.bth.se" into a single run (this also drops the
#    spell-check proofErr wrapper around "Tilljander").
# 2) Split " himself. This will be a Singleton class to be shared over
#    all screens" so the sentence starting at the period is highlighted
#    yellow.
# 3) Split "Object for handling the enemy AI. Moves using a grid
#    structure with nodes." so that just the word "nodes" is
#    highlighted yellow.

$d = $word.ActiveDocument

# --- 1) Merge "Rasmus " + "Tilljander" + " - rati10@student.bth.se" ---
$d.Content.Find.Execute(
    "Rasmus Tilljander - rati10@student.bth.se", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Rasmus Tilljander - rati10@student.bth.se", 2) | Out-Null

# --- 2) Highlight ". This will be a Singleton class to be shared over
#         all screens" (everything from the period onward) yellow ---
$himself = $d.Content
$himself.Find.Execute("himself") | Out-Null
$afterHimself = $himself.End

$singleton = $d.Content
$singleton.Start = $afterHimself
$singletonFind = $singleton.Find
$singletonFind.ClearFormatting()
$singletonFind.Replacement.ClearFormatting()
$singletonFind.Replacement.Highlight = 1
$singletonFind.Execute(
    ". This will be a Singleton class to be shared over all screens",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ". This will be a Singleton class to be shared over all screens", 2) | Out-Null

# --- 3) Highlight just "nodes" in the Ghost paragraph yellow ---
$ghost = $d.Content
$ghost.Find.Execute(
    "Object for handling the enemy AI. Moves using a grid structure with nodes.") | Out-Null
$ghostStart = $ghost.Start

$nodes = $d.Content
$nodes.Start = $ghostStart
$nodesFind = $nodes.Find
$nodesFind.ClearFormatting()
$nodesFind.Replacement.ClearFormatting()
$nodesFind.Replacement.Highlight = 1
$nodesFind.Execute("nodes", $false, $false, $false, $false, $false, $true, 1,
    $false, "nodes", 2) | Out-Null
